$p = $ppt.ActivePresentation

# 1) Update the cached "datetimeFigureOut" field text from "5/31/2024" to
#    "6/1/24" everywhere it appears: the slide master and every slide
#    layout's Date Placeholder shape.
$oldDate = "5/31/2024"
$newDate = "6/1/24"

$m = $p.SlideMaster

for ($j = 1; $j -le $m.Shapes.Count; $j++) {
    $sh = $m.Shapes.Item($j)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

for ($k = 1; $k -le $m.CustomLayouts.Count; $k++) {
    $lay = $m.CustomLayouts.Item($k)
    for ($j = 1; $j -le $lay.Shapes.Count; $j++) {
        $sh = $lay.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# 2) Remove the "TextBox 5" shape (House Budget Committee / cosine
#    similarity caption) from slide 7.
$s7 = $p.Slides.Item(7)
for ($j = $s7.Shapes.Count; $j -ge 1; $j--) {
    $sh = $s7.Shapes.Item($j)
    if ($sh.Name -eq "TextBox 5") {
        $sh.Delete()
    }
}
